$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("B16").Value = 7.12
$ws.Range("C16").Value = -9.09
$ws.Range("E24").Value = 90

# Update the view: scroll so row 13 is the top visible row, and select E24
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("E24").Select()
